$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet holds one row per (market / product / quality) price quote for a
# given week. A new week of data (2023-07-14, serial 45121) is being added
# for "Agricola del Norte S.A. de Arica" - Palta - Hass.
#
# The existing last three rows (224-226, week 2023-05-26 / serial 45072)
# are pushed down to rows 227-229 unchanged, and rows 224-226 are
# overwritten with the new week's figures (which also introduces a new
# "Especial" quality grade ahead of "Primera"/"Segunda").
# ---------------------------------------------------------------------------

# --- Step 1: preserve the old week's rows by writing them to rows 227-229 ---

$oldRows = @(227, 228, 229)
$oldQuality = @("Primera", "Segunda", "Tercera")
$oldMin = @(28000, 25000, 23000)
$oldMax = @(29000, 26000, 24000)
$oldAvg = @(28500, 25500, 23500)
$oldKg = @(2850, 2550, 2350)

for ($i = 0; $i -lt 3; $i++) {
    $r = $oldRows[$i]

    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($r, 4).Value = 45072
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 15
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100106
    $ws.Cells.Item($r, 8).Value = "Oleaginosos"
    $ws.Cells.Item($r, 9).Value = 100106002
    $ws.Cells.Item($r, 10).Value = "Palta"
    $ws.Cells.Item($r, 11).Value = "Hass"
    $ws.Cells.Item($r, 12).Value = $oldQuality[$i]
    $ws.Cells.Item($r, 13).Value = 208
    $ws.Cells.Item($r, 14).Value = $oldMin[$i]
    $ws.Cells.Item($r, 15).Value = $oldMax[$i]
    $ws.Cells.Item($r, 16).Value = $oldAvg[$i]
    $ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = $oldKg[$i]
    $ws.Cells.Item($r, 20).Value = 10
}

# --- Step 2: overwrite rows 224-226 with the new week's figures ---

$newRows = @(224, 225, 226)
$newQuality = @("Especial", "Primera", "Segunda")
$newMin = @(22000, 20000, 18000)
$newMax = @(23000, 21000, 19000)
$newAvg = @(22500, 20500, 18500)
$newKg = @(2250, 2050, 1850)

for ($i = 0; $i -lt 3; $i++) {
    $r = $newRows[$i]

    $ws.Cells.Item($r, 4).Value = 45121
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 12).Value = $newQuality[$i]
    $ws.Cells.Item($r, 14).Value = $newMin[$i]
    $ws.Cells.Item($r, 15).Value = $newMax[$i]
    $ws.Cells.Item($r, 16).Value = $newAvg[$i]
    $ws.Cells.Item($r, 19).Value = $newKg[$i]
}

$wb.Save()
